# Update the NATMI LR-pairs sheet (Col4a3 -> Cd93) with newly recomputed TPM-based values.
# Columns G/H = ligand average/total expression (depends on sending cluster only)
# Columns I/J = ligand derived specificity of avg/total expression
# Columns M/N = receptor average/total expression (depends on target cluster only)
# Columns O/P = receptor derived specificity of avg/total expression
# Columns Q/R = edge average/total expression weight (G*M / H*N)
# Columns S/T = edge derived specificity of avg/total expression weight

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> column -> new value, exactly matching the recomputed TPM results
$updates = @{
    2  = @{ G = 0.6955149999999999; H = 2.086545; I = 0.04699474681458952; J = 0.04699474681458952;
            M = 83.91225566666667;  N = 251.736767; O = 0.9556261553553385; P = 0.9556261553553385;
            Q = 58.36223250000166;  R = 525.2600925000149; S = 0.04490940922032372; T = 0.04490940922032372 }
    3  = @{ G = 0.6955149999999999; H = 2.086545; I = 0.04699474681458952; J = 0.04699474681458952;
            O = 0.00439999103960854; P = 0.00439999103960854;
            Q = 0.2687173206933333; R = 2.41845588624; S = 0.0002067764648928659; T = 0.0002067764648928659 }
    4  = @{ G = 0.6955149999999999; H = 2.086545; I = 0.04699474681458952; J = 0.04699474681458952;
            M = 3.510050666666667;  N = 10.530152; O = 0.03997385360505296; P = 0.03997385360505297;
            Q = 2.441292889426666;  R = 21.97163600484; S = 0.001878561129372931; T = 0.001878561129372931 }
    5  = @{ I = 0.8773676626402658; J = 0.8773676626402658;
            M = 83.91225566666667;  N = 251.736767; O = 0.9556261553553385; P = 0.9556261553553385;
            Q = 1089.592752079634;  R = 9806.334768716706; S = 0.8384354862820168; T = 0.8384354862820168 }
    6  = @{ I = 0.8773676626402658; J = 0.8773676626402658;
            O = 0.00439999103960854; P = 0.00439999103960854;
            S = 0.003860409854059458; T = 0.003860409854059458 }
    7  = @{ I = 0.8773676626402658; J = 0.8773676626402658;
            M = 3.510050666666667;  N = 10.530152; O = 0.03997385360505296; P = 0.03997385360505297;
            Q = 45.577677961904;    R = 410.1991016571361; S = 0.03507176650418948; T = 0.03507176650418949 }
    8  = @{ G = 1.119424666666667;  H = 3.358274; I = 0.07563759054514464; J = 0.07563759054514466;
            M = 83.91225566666667;  N = 251.736767; O = 0.9556261553553385; P = 0.9556261553553385;
            Q = 93.93344882890645;  R = 845.4010394601581; S = 0.07228125985299787; T = 0.07228125985299788 }
    9  = @{ G = 1.119424666666667;  H = 3.358274; I = 0.07563759054514464; J = 0.07563759054514466;
            O = 0.00439999103960854; P = 0.00439999103960854;
            Q = 0.432497929080889;  R = 3.892481361728001; S = 0.0003328047206562161; T = 0.0003328047206562161 }
    10 = @{ G = 1.119424666666667;  H = 3.358274; I = 0.07563759054514464; J = 0.07563759054514466;
            M = 3.510050666666667;  N = 10.530152; O = 0.03997385360505296; P = 0.03997385360505297;
            Q = 3.929237297516445;  R = 35.363135677648; S = 0.00302352597149055; T = 0.003023525971490551 }
}

foreach ($row in $updates.Keys) {
    $rowValues = $updates[$row]
    foreach ($col in $rowValues.Keys) {
        $ws.Range("$col$row").Value = $rowValues[$col]
    }
}
